# Add season-record columns (Wins / Losses / Ties) to the player table.
#
# The sheet previously ended at column AC ("Unnamed: 28"); we extend it with
# three new columns: AD = Wins, AE = Losses, AF = Ties, filled in for every
# player row (2-68) with the team's season record (73 wins, 89 losses, 0 ties).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from the last existing header cell (AC1) onto the
# three new header cells so they match the rest of the header row (bold,
# bordered, centered) instead of creating a brand-new style entry.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record for every data row.
$ws.Range("AD2:AD68").Value = 73
$ws.Range("AE2:AE68").Value = 89
$ws.Range("AF2:AF68").Value = 0
